# Apply PIP_GEN updates to "Parts - Consoles.xlsx"
$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Parts - Consoles" ---
$ws1 = $wb.Worksheets.Item("Parts - Consoles")

# Part Number (D2): 398332 -> 405368 (keep as text, not a number)
$ws1.Range("D2").Value = "'405368"

# Description (E2): ETNT17217V1 -> ETNT24019
$ws1.Range("E2").Value = "ETNT24019"

# Fan Test (AS2): rewritten instructions
$ws1.Range("AS2").Value = '“在控制台上，按大风扇按钮打开风扇，风扇应该低速运行。
再次按大风扇按钮，风扇应该运行在高位。
按小风扇按钮将风扇转到低位，然后再次按下以关闭风扇。“

On the console, press the Large Fan button to turn fan on. Fan should run on low.
Press the Large Fan button again and the fan should run on high.
Press the Small Fan button to turn fan to low, then press again to turn fan off.'

# Updated timestamp (BA2)
$ws1.Range("BA2").Value = 43383.40542136635

# --- Sheet 2: "revision" ---
$ws2 = $wb.Worksheets.Item("revision")

# Revision log timestamp (B14)
$ws2.Range("B14").Value = 43383.40522267423
